$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 33; every row that used to
# be at 33..129 shifts down by one (to 34..130). Insert a blank row at 33
# (this also extends the used range / dimension to R130) then fill it in.
$ws.Rows(33).Insert()

$ws.Cells.Item(33, 1).Value = 2
$ws.Cells.Item(33, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 45274
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 100112026
$ws.Cells.Item(33, 7).Value = "Haba"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 11000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(33, 16).Value = 440
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
